$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 284.25
$ws.Range("I41").Value = 215.85715
$ws.Range("J41").Value = 380
$ws.Range("K41").Value = 215.85715
$ws.Range("L41").Value = 380
$ws.Range("M41").Value = 224.14285
$ws.Range("N41").Value = -1260
$ws.Range("H74").Value = 4850.077
$ws.Range("I74").Value = 4364.7144
$ws.Range("J74").Value = 5416.3335
$ws.Range("K74").Value = 4364.7144
$ws.Range("L74").Value = 5416.3335
$ws.Range("M74").Value = -3428.7144
$ws.Range("N74").Value = -7288.3335
$ws.Range("H77").Value = 4850.077
$ws.Range("I77").Value = 4364.7144
$ws.Range("J77").Value = 5416.3335
$ws.Range("K77").Value = 21823.572
$ws.Range("L77").Value = 27081.6675
$ws.Range("M77").Value = -17143.572
$ws.Range("N77").Value = -36441.6675
$ws.Range("H80").Value = 6068.15
$ws.Range("I80").Value = 513.8182
$ws.Range("J80").Value = 12856.777
$ws.Range("K80").Value = 1541.4546
$ws.Range("L80").Value = 38570.331
$ws.Range("M80").Value = -543.4546
$ws.Range("N80").Value = -40566.331
$ws.Range("H83").Value = 6068.15
$ws.Range("I83").Value = 513.8182
$ws.Range("J83").Value = 12856.777
$ws.Range("K83").Value = 4624.3638
$ws.Range("L83").Value = 115710.993
$ws.Range("M83").Value = 367.6361999999999
$ws.Range("N83").Value = -125694.993
$ws.Range("H94").Value = 2800
$ws.Range("I94").Value = 2800
$ws.Range("K94").Value = 2800
$ws.Range("M94").Value = -2349
$ws.Range("H100").Value = 1651.9
$ws.Range("I100").Value = 1601.875
$ws.Range("J100").Value = 1852
$ws.Range("K100").Value = 1601.875
$ws.Range("L100").Value = 1852
$ws.Range("M100").Value = -1060.875
$ws.Range("N100").Value = -2934
$ws.Range("H103").Value = 1325
$ws.Range("J103").Value = 1341.6666
$ws.Range("L103").Value = 4024.9998
$ws.Range("N103").Value = -5196.9998
$ws.Range("H116").Value = 2143.0356
$ws.Range("I116").Value = 1965
$ws.Range("J116").Value = 2321.0715
$ws.Range("K116").Value = 1965
$ws.Range("L116").Value = 2321.0715
$ws.Range("M116").Value = 1477
$ws.Range("N116").Value = -9205.0715
$ws.Range("H128").Value = 39271.43
$ws.Range("J128").Value = 39271.43
$ws.Range("L128").Value = 39271.43
$ws.Range("N128").Value = -49231.43
$ws.Range("H129").Value = 919.1948
$ws.Range("I129").Value = 370.57144
$ws.Range("J129").Value = 974.0571
$ws.Range("K129").Value = 1111.71432
$ws.Range("L129").Value = 2922.1713
$ws.Range("M129").Value = 3888.28568
$ws.Range("N129").Value = -12922.1713
$ws.Range("H137").Value = 3952.8096
$ws.Range("I137").Value = 6675.75
$ws.Range("K137").Value = 20027.25
$ws.Range("M137").Value = -17477.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 20000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 20000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 20000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -20300
$ws.Range("H41").Value = 8739
$ws.Range("I41").Value = 4956
$ws.Range("K41").Value = 4956
$ws.Range("M41").Value = -4542
$ws.Range("H45").Value = 1549.7084
$ws.Range("I45").Value = 1577.9565
$ws.Range("J45").Value = 900
$ws.Range("K45").Value = 1577.9565
$ws.Range("L45").Value = 900
$ws.Range("M45").Value = -1200.9565
$ws.Range("N45").Value = -1654
$ws.Range("H88").Value = 6834
$ws.Range("I88").Value = 18160
$ws.Range("J88").Value = 3058.6667
$ws.Range("K88").Value = 18160
$ws.Range("L88").Value = 3058.6667
$ws.Range("M88").Value = -17754
$ws.Range("N88").Value = -3870.6667
$ws.Range("H91").Value = 6834
$ws.Range("I91").Value = 18160
$ws.Range("J91").Value = 3058.6667
$ws.Range("K91").Value = 18160
$ws.Range("L91").Value = 3058.6667
$ws.Range("M91").Value = -16756
$ws.Range("N91").Value = -5866.6667
$ws.Range("H97").Value = 1142.5834
$ws.Range("I97").Value = 1150
$ws.Range("K97").Value = 1150
$ws.Range("M97").Value = -654
$ws.Range("H110").Value = 1102.1154
$ws.Range("I110").Value = 1067.6086
$ws.Range("J110").Value = 1366.6666
$ws.Range("K110").Value = 1067.6086
$ws.Range("L110").Value = 1366.6666
$ws.Range("M110").Value = 977.3914
$ws.Range("N110").Value = -5456.6666
$ws.Range("H122").Value = 11365724
$ws.Range("I122").Value = 2260.7778
$ws.Range("J122").Value = 62501308
$ws.Range("K122").Value = 6782.3334
$ws.Range("L122").Value = 187503924
$ws.Range("M122").Value = -4332.3334
$ws.Range("N122").Value = -187508824
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 627.7727
$ws.Range("I94").Value = 705
$ws.Range("J94").Value = 492.625
$ws.Range("K94").Value = 705
$ws.Range("L94").Value = 492.625
$ws.Range("M94").Value = -254
$ws.Range("N94").Value = -1394.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7972.4443
$ws.Range("I31").Value = 6148
$ws.Range("J31").Value = 16000
$ws.Range("K31").Value = 6148
$ws.Range("L31").Value = 16000
$ws.Range("M31").Value = -5853
$ws.Range("N31").Value = -16590
$ws.Range("H34").Value = 7972.4443
$ws.Range("I34").Value = 6148
$ws.Range("J34").Value = 16000
$ws.Range("K34").Value = 6148
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = -5946
$ws.Range("N34").Value = -16404
$ws.Range("H121").Value = 24000
$ws.Range("J121").Value = 24000
$ws.Range("L121").Value = 24000
$ws.Range("N121").Value = -26620
$ws.Range("H127").Value = 3808895
$ws.Range("J127").Value = 3808895
$ws.Range("L127").Value = 3808895
$ws.Range("N127").Value = -3818815
$ws.Range("H132").Value = 3485
$ws.Range("I132").Value = 2707.923
$ws.Range("J132").Value = 5168.6665
$ws.Range("K132").Value = 8123.768999999999
$ws.Range("L132").Value = 15505.9995
$ws.Range("M132").Value = -5593.768999999999
$ws.Range("N132").Value = -20565.9995
$ws.Range("H133").Value = 41997.5
$ws.Range("J133").Value = 41997.5
$ws.Range("L133").Value = 41997.5
$ws.Range("N133").Value = -47057.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8107.5713
$ws.Range("I68").Value = 637.5
$ws.Range("J68").Value = 11095.6
$ws.Range("K68").Value = 1912.5
$ws.Range("L68").Value = 33286.8
$ws.Range("M68").Value = -1101.5
$ws.Range("N68").Value = -34908.8
$ws.Range("H71").Value = 8107.5713
$ws.Range("I71").Value = 637.5
$ws.Range("J71").Value = 11095.6
$ws.Range("K71").Value = 5737.5
$ws.Range("L71").Value = 99860.40000000001
$ws.Range("M71").Value = -1681.5
$ws.Range("N71").Value = -107972.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1254.1666
$ws.Range("J97").Value = 1050.25
$ws.Range("L97").Value = 1050.25
$ws.Range("N97").Value = -2042.25
$ws.Range("H126").Value = 1735.5686
$ws.Range("J126").Value = 1722.7
$ws.Range("L126").Value = 5168.1
$ws.Range("N126").Value = -10108.1
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6083.1665
$ws.Range("I7").Value = 6083.1665
$ws.Range("K7").Value = 6083.1665
$ws.Range("M7").Value = -5971.1665
$ws.Range("H16").Value = 1001
$ws.Range("I16").Value = 1001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -831
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 3466.625
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 4493.2
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 4493.2
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -4765.2
$ws.Range("H93").Value = 1167.1818
$ws.Range("I93").Value = 945
$ws.Range("J93").Value = 1433.8
$ws.Range("K93").Value = 945
$ws.Range("L93").Value = 1433.8
$ws.Range("M93").Value = 303
$ws.Range("N93").Value = -3929.8
$ws.Range("H100").Value = 5191.8
$ws.Range("I100").Value = 4635
$ws.Range("K100").Value = 4635
$ws.Range("M100").Value = -4094
$ws.Range("H126").Value = 6083.1665
$ws.Range("I126").Value = 6083.1665
$ws.Range("K126").Value = 18249.4995
$ws.Range("M126").Value = -15779.4995
$ws.Range("H129").Value = 69482.25
$ws.Range("J129").Value = 69482.25
$ws.Range("L129").Value = 69482.25
$ws.Range("N129").Value = -79482.25
$ws.Range("H139").Value = 54805
$ws.Range("J139").Value = 54805
$ws.Range("L139").Value = 54805
$ws.Range("N139").Value = -65085
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14289234
$ws.Range("I81").Value = 1880
$ws.Range("J81").Value = 15388261
$ws.Range("K81").Value = 3760
$ws.Range("L81").Value = 30776522
$ws.Range("M81").Value = -2699
$ws.Range("N81").Value = -30778644
$ws.Range("H84").Value = 14289234
$ws.Range("I84").Value = 1880
$ws.Range("J84").Value = 15388261
$ws.Range("K84").Value = 18800
$ws.Range("L84").Value = 153882610
$ws.Range("M84").Value = -13496
$ws.Range("N84").Value = -153893218
$ws.Range("H122").Value = 6150.7617
$ws.Range("I122").Value = 1580.5454
$ws.Range("J122").Value = 11178
$ws.Range("K122").Value = 4741.6362
$ws.Range("L122").Value = 33534
$ws.Range("M122").Value = -2291.6362
$ws.Range("N122").Value = -38434
$ws.Range("H132").Value = 3380.65
$ws.Range("I132").Value = 3163.4375
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 9490.3125
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -6960.3125
$ws.Range("N132").Value = -17808.5
